$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new time-tracking entries (rows 40 and 41) ---
# Set values first so the running-total formula (F4 = SUM(B2:B1000))
# recalculates against the new hours before we touch formatting.
$ws.Range("A40").Value = 45639
$ws.Range("B40").Value = 1.5
$ws.Range("C40").Value = "Review callstack and architecture resources"

$ws.Range("A41").Value = 45639
$ws.Range("B41").Value = 1.5
$ws.Range("C41").Value = "Review assembly code resources"

# Copy formatting from the previous row (row 39) onto the new rows so the
# date column keeps its date number format / styling, matching the rest
# of the table.
$ws.Range("A39").Copy()
$ws.Range("A40:A41").PasteSpecial(-4122)
$ws.Range("B39").Copy()
$ws.Range("B40:B41").PasteSpecial(-4122)
$ws.Range("C39").Copy()
$ws.Range("C40:C41").PasteSpecial(-4122)

# --- Update the view state to reflect the newly added rows ---
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("A42").Select() | Out-Null

# Make sure everything (the running total formula) is recalculated.
$excel.CalculateFull()
